# "updated to latest version"
#
# Product Backlog.xlsx: a new backlog item's Priority/Complexity scores
# were filled in, the active selection moved, and the workbook window
# was resized on the author's machine.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 ("As a business, we want to have all of our functionality intact
# when implementing student code") gets its Priority (B10) and
# Complexity (C10) scores filled in.
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = 4

# The author's cursor/selection ends up on A7 instead of A22.
$ws.Range("A7").Select()

# The workbook window was resized (from 21570x8055 to 17190x7260 twips)
# when the file was last saved.
$excel.ActiveWindow.Width = 17190
$excel.ActiveWindow.Height = 7260
